$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = $null
$ws.Range("C2").Value = 0.68758068472521283
$ws.Range("D2").Value = 0.38194044403600752
$ws.Range("E2").Value = 1.0810840941366469

$ws.Range("B3").Value = 1.5369245927528954
$ws.Range("C3").Value = 1.9833558962570397
$ws.Range("D3").Value = 0.88281367201455507
$ws.Range("E3").Value = 1.3291069666644408

$ws.Range("B1:E3").Select()
